# Workbook: dataEDT-ELT-S2-2026.xlsx
# - "Feuille2" worksheet holds the full (filtered) schedule table (A1:G559).
# - Correct a missing room value for the HADER / Cours-Technologie... row (F261: "-" -> "S14").
# - Add an AutoFilter criterion on the "Enseignants" column (C, colId=2) for "HADER",
#   layered on top of the pre-existing "Promotion" (G, colId=6) "ING2" filter.
# - Update the saved cursor/selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the room ("Lieu") value that was left blank ("-") for the HADER / TDEE row.
$ws.Range("F261").Value = "S14"

# Re-apply the existing filter (Promotion = ING2, field 7) together with the new
# filter (Enseignants = HADER, field 3) so the new filterColumn is emitted ahead
# of the pre-existing one, matching how the sheet was re-saved.
$ws.Range("A1:G559").AutoFilter(3, "HADER", 7)
$ws.Range("A1:G559").AutoFilter(7, "ING2", 7)

# Move the saved selection/active cell to F566 (below the data, as left by the editor).
$ws.Range("F566").Select()
